# Append three new match rows (118-120) to the Ekstraklasa 2023-2024 sheet,
# mirroring the formatting of the last existing data row (117).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data, in column order A..V (A is recomputed as row-1 below).
$newRows = @(
    @{
        F="Puszcza"; G=0; H="Pogon Szczecin"; I=2;
        E=45234.625;
        J=3.83;  K="28/10/2023 18:13"; L=4.68; M="04/11/2023 14:32";
        N=3.69;  O="28/10/2023 18:13"; P=4.13; Q="04/11/2023 14:44";
        R=1.95;  S="28/10/2023 18:13"; T=1.72; U="04/11/2023 14:32";
        V="https://www.betexplorer.com/football/poland/ekstraklasa/puszcza-pogon-szczecin/rX9PLxGT/"
    },
    @{
        F="Widzew Lodz"; G=0; H="Warta Poznan"; I=1;
        E=45234.72916666666;
        J=2.3;   K="29/10/2023 17:43"; L=2.33; M="04/11/2023 17:29";
        N=3.18;  O="29/10/2023 17:43"; P=3.06; Q="04/11/2023 17:29";
        R=3.4;   S="29/10/2023 17:43"; T=3.59; U="04/11/2023 17:10";
        V="https://www.betexplorer.com/football/poland/ekstraklasa/widzew-lodz-warta-poznan/Cnlwu0Gi/"
    },
    @{
        F="Lech Poznan"; G=2; H="Ruch Chorzow"; I=0;
        E=45234.83333333334;
        J=1.42;  K="28/10/2023 21:12"; L=1.36; M="04/11/2023 19:51";
        N=4.76;  O="28/10/2023 21:12"; P=5.21; Q="04/11/2023 19:55";
        R=7.54;  S="28/10/2023 21:12"; T=8.94; U="04/11/2023 19:55";
        V="https://www.betexplorer.com/football/poland/ekstraklasa/lech-poznan-ruch-chorzow/OUgONboH/"
    }
)

$lastRow = 117

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $targetRow = $lastRow + 1 + $i
    $srcRange = $ws.Range("A$lastRow`:V$lastRow")
    $dstRange = $ws.Range("A$targetRow`:V$targetRow")
    $srcRange.Copy($dstRange)

    $data = $newRows[$i]

    $ws.Cells.Item($targetRow, 1).Value = $targetRow - 1   # A: Indice
    $ws.Cells.Item($targetRow, 2).Value = "poland"         # B: pais
    $ws.Cells.Item($targetRow, 3).Value = "ekstraklasa"    # C: torneio
    $ws.Cells.Item($targetRow, 4).Value = "2023-2024"      # D: temporada
    $ws.Cells.Item($targetRow, 5).Value = $data.E          # E: data_partida
    $ws.Cells.Item($targetRow, 6).Value = $data.F          # F: home
    $ws.Cells.Item($targetRow, 7).Value = $data.G          # G: home_ft_gols
    $ws.Cells.Item($targetRow, 8).Value = $data.H          # H: away
    $ws.Cells.Item($targetRow, 9).Value = $data.I          # I: away_ft_gols
    $ws.Cells.Item($targetRow, 10).Value = $data.J         # J: home_opening_odds
    $ws.Cells.Item($targetRow, 11).Value = $data.K         # K: home_opening_data_hora
    $ws.Cells.Item($targetRow, 12).Value = $data.L         # L: home_closing_odds
    $ws.Cells.Item($targetRow, 13).Value = $data.M         # M: home_closing_data_hora
    $ws.Cells.Item($targetRow, 14).Value = $data.N         # N: draw_opening_odds
    $ws.Cells.Item($targetRow, 15).Value = $data.O         # O: draw_opening_data_hora
    $ws.Cells.Item($targetRow, 16).Value = $data.P         # P: draw_closing_odds
    $ws.Cells.Item($targetRow, 17).Value = $data.Q         # Q: draw_closing_data_hora
    $ws.Cells.Item($targetRow, 18).Value = $data.R         # R: away_opening_odds
    $ws.Cells.Item($targetRow, 19).Value = $data.S         # S: away_opening_data_hora
    $ws.Cells.Item($targetRow, 20).Value = $data.T         # T: away_closing_odds
    $ws.Cells.Item($targetRow, 21).Value = $data.U         # U: away_closing_data_hora
    $ws.Cells.Item($targetRow, 22).Value = $data.V         # V: url_partida
}
